$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated benchmark numbers (K/L columns); the M column speedup
#     formulas are already in the sheet and recalculate automatically. ---
$ws.Range("K9").Value = 11
$ws.Range("L9").Value = 26

$ws.Range("K10").Value = 17
$ws.Range("L10").Value = 25

$ws.Range("K11").Value = 86
$ws.Range("L11").Value = 122

$ws.Range("K12").Value = 623
$ws.Range("L12").Value = 966

$ws.Range("K13").Value = 5725
$ws.Range("L13").Value = 12097

# --- New (empty) formatted cell a few rows below the table, as if a
#     console-style font was applied in preparation for more data. ---
$consoleCell = $ws.Range("L19")
$consoleCell.Font.Family = 3
$consoleCell.Font.Name = "Lucida Console"
$consoleCell.Font.Size = 10
$consoleCell.Font.Color = 0
$consoleCell.VerticalAlignment = -4108

# --- Update the view/selection state to match where the user ended up. ---
[void]$ws.Range("R20").Select()
